# Apply the recorded edits to the "Hoja1" worksheet (the sheet with
# dimension A1:N13 in the original workbook).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Activate()

# New entries in column H for several expense rows.
$ws.Range("H2").Value = 4
$ws.Range("H3").Value = 4
$ws.Range("H4").Value = 4
$ws.Range("H6").Value = 4
$ws.Range("H8").Value = 44

# Totals (N column) and the grand-total row (13) recalculate automatically
# from the SUM formulas already present in the sheet.

# Leave the selection/active cell on H12, matching the saved view state.
$null = $ws.Range("H12").Select()
